$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $style = $p.Style.NameLocal
    if ($style -ne "Heading 2") {
        continue
    }
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)

    if ($text -eq "Vista de context") {
        # Insert a new run containing "o" right before the paragraph mark
        # (after the existing bookmarkEnd), so "Vista de context" becomes
        # "Vista de contexto", matching pt-PT language formatting.
        $r = $p.Range
        $insertRange = $d.Range($r.End - 1, $r.End - 1)
        $insertRange.InsertAfter("o")
        $insertRange.Font.LanguageID = "pt-PT"
    }
    elseif ($text -eq "Vista functional") {
        # Fix misspelling "functional" -> "funcional"
        $r = $p.Range
        $r.Find.Execute("functional", $true, $false, $false, $false, $false,
                         $true, 1, $false, "funcional", 2)
    }
}
